$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete rows (old rows 8-10, the third "MuSCs -> *" block)
$ws.Rows("8:10").Delete()

# Recomputed TPM-based values now populate rows 2-7 (target cluster "ECs" rows dropped)
$data = New-Object 'object[,]' 6,20
$data[0,0] = "ECs"
$data[0,1] = "Efna2"
$data[0,2] = "Epha3"
$data[0,3] = "FAPs"
$data[0,4] = 2
$data[0,5] = 0.6666666666666666
$data[0,6] = 0.06446533333333333
$data[0,7] = 0.193396
$data[0,8] = 0.02693738696927793
$data[0,9] = 0.02693738696927793
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 26.097779
$data[0,13] = 78.29333700000001
$data[0,14] = 0.9922055808976035
$data[0,15] = 0.9922055808976036
$data[0,16] = 1.682402022494667
$data[0,17] = 15.141618202452
$data[0,18] = 0.02672742568571594
$data[0,19] = 0.02672742568571594
$data[1,0] = "ECs"
$data[1,1] = "Efna2"
$data[1,2] = "Epha3"
$data[1,3] = "MuSCs"
$data[1,4] = 2
$data[1,5] = 0.6666666666666666
$data[1,6] = 0.06446533333333333
$data[1,7] = 0.193396
$data[1,8] = 0.02693738696927793
$data[1,9] = 0.02693738696927793
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 0.205015
$data[1,13] = 0.6150450000000001
$data[1,14] = 0.007794419102396499
$data[1,15] = 0.007794419102396499
$data[1,16] = 0.01321636031333333
$data[1,17] = 0.11894724282
$data[1,18] = 0.0002099612835619864
$data[1,19] = 0.0002099612835619864
$data[2,0] = "FAPs"
$data[2,1] = "Efna2"
$data[2,2] = "Epha3"
$data[2,3] = "FAPs"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 1.843761666666667
$data[2,7] = 5.531285
$data[2,8] = 0.7704314695358874
$data[2,9] = 0.7704314695358874
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 26.097779
$data[2,13] = 78.29333700000001
$data[2,14] = 0.9922055808976035
$data[2,15] = 0.9922055808976036
$data[2,16] = 48.11808450533834
$data[2,17] = 433.0627605480451
$data[2,18] = 0.7644264037726494
$data[2,19] = 0.7644264037726495
$data[3,0] = "FAPs"
$data[3,1] = "Efna2"
$data[3,2] = "Epha3"
$data[3,3] = "MuSCs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 1.843761666666667
$data[3,7] = 5.531285
$data[3,8] = 0.7704314695358874
$data[3,9] = 0.7704314695358874
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 0.205015
$data[3,13] = 0.6150450000000001
$data[3,14] = 0.007794419102396499
$data[3,15] = 0.007794419102396499
$data[3,16] = 0.3779987980916668
$data[3,17] = 3.401989182825001
$data[3,18] = 0.006005065763237926
$data[3,19] = 0.006005065763237926
$data[4,0] = "MuSCs"
$data[4,1] = "Efna2"
$data[4,2] = "Epha3"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 0.4849276666666666
$data[4,7] = 1.454783
$data[4,8] = 0.2026311434948347
$data[4,9] = 0.2026311434948347
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 26.097779
$data[4,13] = 78.29333700000001
$data[4,14] = 0.9922055808976035
$data[4,15] = 0.9922055808976036
$data[4,16] = 12.65553507565233
$data[4,17] = 113.899815680871
$data[4,18] = 0.2010517514392381
$data[4,19] = 0.2010517514392381
$data[5,0] = "MuSCs"
$data[5,1] = "Efna2"
$data[5,2] = "Epha3"
$data[5,3] = "MuSCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 0.4849276666666666
$data[5,7] = 1.454783
$data[5,8] = 0.2026311434948347
$data[5,9] = 0.2026311434948347
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 0.205015
$data[5,13] = 0.6150450000000001
$data[5,14] = 0.007794419102396499
$data[5,15] = 0.007794419102396499
$data[5,16] = 0.09941744558166668
$data[5,17] = 0.8947570102350001
$data[5,18] = 0.001579392055596585
$data[5,19] = 0.001579392055596585

$ws.Range("A2:T7").Value2 = $data
